# Auto-generated Excel COM-interop script
# Applies numeric corrections to H/I/J/K/L/M/N columns across multiple rows
# in sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, matching the target diff.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 462.6129   # H17: 515.13635 -> 462.6129
$ws.Cells.Item(17, 10).Value = 462.6129   # J17: 515.13635 -> 462.6129
$ws.Cells.Item(17, 12).Value = 1387.8387   # L17: 1545.40905 -> 1387.8387
$ws.Cells.Item(17, 14).Value = -1723.8387   # N17: -1881.40905 -> -1723.8387
$ws.Cells.Item(53, 8).Value = 624.1111   # H53: 634.2222 -> 624.1111
$ws.Cells.Item(53, 9).Value = 318.75   # I53: 348.5 -> 318.75
$ws.Cells.Item(53, 10).Value = 868.4   # J53: 862.8 -> 868.4
$ws.Cells.Item(53, 11).Value = 318.75   # K53: 348.5 -> 318.75
$ws.Cells.Item(53, 12).Value = 868.4   # L53: 862.8 -> 868.4
$ws.Cells.Item(53, 13).Value = 318.25   # M53: 288.5 -> 318.25
$ws.Cells.Item(53, 14).Value = -2142.4   # N53: -2136.8 -> -2142.4
$ws.Cells.Item(100, 8).Value = 6186.75   # H100: 7031.4287 -> 6186.75
$ws.Cells.Item(100, 9).Value = 6186.75   # I100: 7031.4287 -> 6186.75
$ws.Cells.Item(100, 11).Value = 6186.75   # K100: 7031.4287 -> 6186.75
$ws.Cells.Item(100, 13).Value = -5645.75   # M100: -6490.4287 -> -5645.75
$ws.Cells.Item(137, 8).Value = 8849.091   # H137: 12205.857 -> 8849.091
$ws.Cells.Item(137, 9).Value = 2974.5   # I137: 3000 -> 2974.5
$ws.Cells.Item(137, 10).Value = 12206   # J137: 13740.167 -> 12206
$ws.Cells.Item(137, 11).Value = 8923.5   # K137: 9000 -> 8923.5
$ws.Cells.Item(137, 12).Value = 36618   # L137: 41220.501 -> 36618
$ws.Cells.Item(137, 13).Value = -6373.5   # M137: -6450 -> -6373.5
$ws.Cells.Item(137, 14).Value = -41718   # N137: -46320.501 -> -41718

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1539.5294   # H2: 1485.1666 -> 1539.5294
$ws.Cells.Item(2, 9).Value = 1521.5333   # I2: 1461.5 -> 1521.5333
$ws.Cells.Item(2, 11).Value = 1521.5333   # K2: 1461.5 -> 1521.5333
$ws.Cells.Item(2, 13).Value = -1408.5333   # M2: -1348.5 -> -1408.5333
$ws.Cells.Item(27, 8).Value = 10000   # H27: 0 -> 10000
$ws.Cells.Item(27, 10).Value = 10000   # J27: 0 -> 10000
$ws.Cells.Item(27, 12).Value = 10000   # L27: 0 -> 10000
$ws.Cells.Item(27, 14).Value = -10368   # N27: None -> -10368
$ws.Cells.Item(32, 8).Value = 9435235   # H32: 9805242 -> 9435235
$ws.Cells.Item(32, 9).Value = 9435235   # I32: 9805242 -> 9435235
$ws.Cells.Item(32, 11).Value = 9435235   # K32: 9805242 -> 9435235
$ws.Cells.Item(32, 13).Value = -9434948   # M32: -9804955 -> -9434948
$ws.Cells.Item(38, 8).Value = 1019   # H38: 0 -> 1019
$ws.Cells.Item(38, 9).Value = 1019   # I38: 0 -> 1019
$ws.Cells.Item(38, 11).Value = 1019   # K38: 0 -> 1019
$ws.Cells.Item(38, 13).Value = -552   # M38: None -> -552
$ws.Cells.Item(45, 8).Value = 2211.3635   # H45: 2082.7917 -> 2211.3635
$ws.Cells.Item(45, 9).Value = 1850.0667   # I45: 1863.3334 -> 1850.0667
$ws.Cells.Item(45, 10).Value = 2985.5715   # J45: 2448.5557 -> 2985.5715
$ws.Cells.Item(45, 11).Value = 1850.0667   # K45: 1863.3334 -> 1850.0667
$ws.Cells.Item(45, 12).Value = 2985.5715   # L45: 2448.5557 -> 2985.5715
$ws.Cells.Item(45, 13).Value = -1473.0667   # M45: -1486.3334 -> -1473.0667
$ws.Cells.Item(45, 14).Value = -3739.5715   # N45: -3202.5557 -> -3739.5715
$ws.Cells.Item(52, 8).Value = 119966   # H52: 119989.5 -> 119966
$ws.Cells.Item(52, 10).Value = 119966   # J52: 119989.5 -> 119966
$ws.Cells.Item(52, 12).Value = 119966   # L52: 119989.5 -> 119966
$ws.Cells.Item(52, 14).Value = -120602   # N52: -120625.5 -> -120602
$ws.Cells.Item(61, 8).Value = 26375986   # H61: 25057288 -> 26375986
$ws.Cells.Item(61, 9).Value = 38467444   # I61: 35719910 -> 38467444
$ws.Cells.Item(61, 11).Value = 38467444   # K61: 35719910 -> 38467444
$ws.Cells.Item(61, 13).Value = -38467232   # M61: -35719698 -> -38467232
$ws.Cells.Item(74, 8).Value = 14717719   # H74: 11914473 -> 14717719
$ws.Cells.Item(74, 9).Value = 25000676   # I74: 17857818 -> 25000676
$ws.Cells.Item(74, 11).Value = 25000676   # K74: 17857818 -> 25000676
$ws.Cells.Item(74, 13).Value = -24999802   # M74: -17856944 -> -24999802
$ws.Cells.Item(77, 8).Value = 14717719   # H77: 11914473 -> 14717719
$ws.Cells.Item(77, 9).Value = 25000676   # I77: 17857818 -> 25000676
$ws.Cells.Item(77, 11).Value = 125003380   # K77: 89289090 -> 125003380
$ws.Cells.Item(77, 13).Value = -124999012   # M77: -89284722 -> -124999012
$ws.Cells.Item(103, 8).Value = 66663.336   # H103: 49250 -> 66663.336
$ws.Cells.Item(103, 10).Value = 66663.336   # J103: 49250 -> 66663.336
$ws.Cells.Item(103, 12).Value = 66663.336   # L103: 49250 -> 66663.336
$ws.Cells.Item(103, 14).Value = -69007.336   # N103: -51594 -> -69007.336
$ws.Cells.Item(116, 8).Value = 1539.5294   # H116: 1485.1666 -> 1539.5294
$ws.Cells.Item(116, 9).Value = 1521.5333   # I116: 1461.5 -> 1521.5333
$ws.Cells.Item(116, 11).Value = 1521.5333   # K116: 1461.5 -> 1521.5333
$ws.Cells.Item(116, 13).Value = 772.4666999999999   # M116: 832.5 -> 772.4666999999999
$ws.Cells.Item(132, 8).Value = 6368.129   # H132: 5631.9165 -> 6368.129
$ws.Cells.Item(132, 9).Value = 4925.3335   # I132: 3674.5 -> 4925.3335
$ws.Cells.Item(132, 10).Value = 7279.3687   # J132: 7589.3335 -> 7279.3687
$ws.Cells.Item(132, 11).Value = 14776.0005   # K132: 11023.5 -> 14776.0005
$ws.Cells.Item(132, 12).Value = 21838.1061   # L132: 22768.0005 -> 21838.1061
$ws.Cells.Item(132, 13).Value = -12246.0005   # M132: -8493.5 -> -12246.0005
$ws.Cells.Item(132, 14).Value = -26898.1061   # N132: -27828.0005 -> -26898.1061
$ws.Cells.Item(136, 8).Value = 26375986   # H136: 25057288 -> 26375986
$ws.Cells.Item(136, 9).Value = 38467444   # I136: 35719910 -> 38467444
$ws.Cells.Item(136, 11).Value = 115402332   # K136: 107159730 -> 115402332
$ws.Cells.Item(136, 13).Value = -115399782   # M136: -107157180 -> -115399782

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1539.5294   # H3: 1485.1666 -> 1539.5294
$ws.Cells.Item(3, 9).Value = 1521.5333   # I3: 1461.5 -> 1521.5333
$ws.Cells.Item(3, 11).Value = 1521.5333   # K3: 1461.5 -> 1521.5333
$ws.Cells.Item(3, 13).Value = -1407.5333   # M3: -1347.5 -> -1407.5333
$ws.Cells.Item(107, 8).Value = 1161.3334   # H107: 1228.7 -> 1161.3334
$ws.Cells.Item(107, 9).Value = 1161.3334   # I107: 1228.7 -> 1161.3334
$ws.Cells.Item(107, 11).Value = 1161.3334   # K107: 1228.7 -> 1161.3334
$ws.Cells.Item(107, 13).Value = 758.6666   # M107: 691.3 -> 758.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4932   # H58: 0 -> 4932
$ws.Cells.Item(58, 9).Value = 4932   # I58: 0 -> 4932
$ws.Cells.Item(58, 11).Value = 4932   # K58: 0 -> 4932
$ws.Cells.Item(58, 13).Value = -4729   # M58: None -> -4729
$ws.Cells.Item(132, 8).Value = 3315.6667   # H132: 3333 -> 3315.6667
$ws.Cells.Item(132, 9).Value = 2827.7144   # I132: 2666.3333 -> 2827.7144
$ws.Cells.Item(132, 10).Value = 3998.8   # J132: 3999.6667 -> 3998.8
$ws.Cells.Item(132, 11).Value = 8483.143199999999   # K132: 7998.999899999999 -> 8483.143199999999
$ws.Cells.Item(132, 12).Value = 11996.4   # L132: 11999.0001 -> 11996.4
$ws.Cells.Item(132, 13).Value = -5953.143199999999   # M132: -5468.999899999999 -> -5953.143199999999
$ws.Cells.Item(132, 14).Value = -17056.4   # N132: -17059.0001 -> -17056.4
$ws.Cells.Item(134, 8).Value = 593206.6   # H134: 720069.0600000001 -> 593206.6
$ws.Cells.Item(134, 9).Value = 910413.4399999999   # I134: 1251375.2 -> 910413.4399999999
$ws.Cells.Item(134, 11).Value = 2731240.32   # K134: 3754125.6 -> 2731240.32
$ws.Cells.Item(134, 13).Value = -2728705.32   # M134: -3751590.6 -> -2728705.32
$ws.Cells.Item(136, 8).Value = 4932   # H136: 0 -> 4932
$ws.Cells.Item(136, 9).Value = 4932   # I136: 0 -> 4932
$ws.Cells.Item(136, 11).Value = 14796   # K136: 0 -> 14796
$ws.Cells.Item(136, 13).Value = -12246   # M136: None -> -12246

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(35, 8).Value = 2200   # H35: 2166.6667 -> 2200
$ws.Cells.Item(35, 9).Value = 2200   # I35: 2166.6667 -> 2200
$ws.Cells.Item(35, 11).Value = 6600   # K35: 6500.000100000001 -> 6600
$ws.Cells.Item(35, 13).Value = -6312   # M35: -6212.000100000001 -> -6312
$ws.Cells.Item(37, 8).Value = 82999.5   # H37: 83994.75 -> 82999.5
$ws.Cells.Item(37, 10).Value = 82999.5   # J37: 83994.75 -> 82999.5
$ws.Cells.Item(37, 12).Value = 248998.5   # L37: 251984.25 -> 248998.5
$ws.Cells.Item(37, 14).Value = -249222.5   # N37: -252208.25 -> -249222.5
$ws.Cells.Item(107, 8).Value = 430.29413   # H107: 445.3125 -> 430.29413
$ws.Cells.Item(107, 9).Value = 424.8889   # I107: 454.25 -> 424.8889
$ws.Cells.Item(107, 11).Value = 1274.6667   # K107: 1362.75 -> 1274.6667
$ws.Cells.Item(107, 13).Value = 645.3333   # M107: 557.25 -> 645.3333
$ws.Cells.Item(131, 8).Value = 5059.316   # H131: 5485.0586 -> 5059.316
$ws.Cells.Item(131, 10).Value = 4022.25   # J131: 4882.8335 -> 4022.25
$ws.Cells.Item(131, 12).Value = 12066.75   # L131: 14648.5005 -> 12066.75
$ws.Cells.Item(131, 14).Value = -22146.75   # N131: -24728.5005 -> -22146.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3206.1707   # H102: 3421.923 -> 3206.1707
$ws.Cells.Item(102, 9).Value = 2708.3872   # I102: 2857.3103 -> 2708.3872
$ws.Cells.Item(102, 10).Value = 4749.3   # J102: 5059.3 -> 4749.3
$ws.Cells.Item(102, 11).Value = 2708.3872   # K102: 2857.3103 -> 2708.3872
$ws.Cells.Item(102, 12).Value = 4749.3   # L102: 5059.3 -> 4749.3
$ws.Cells.Item(102, 13).Value = -1086.3872   # M102: -1235.3103 -> -1086.3872
$ws.Cells.Item(102, 14).Value = -7993.3   # N102: -8303.299999999999 -> -7993.3
$ws.Cells.Item(109, 8).Value = 46774.25   # H109: 46774.5 -> 46774.25
$ws.Cells.Item(109, 10).Value = 46774.25   # J109: 46774.5 -> 46774.25
$ws.Cells.Item(109, 12).Value = 46774.25   # L109: 46774.5 -> 46774.25
$ws.Cells.Item(109, 14).Value = -48854.25   # N109: -48854.5 -> -48854.25
$ws.Cells.Item(122, 8).Value = 1996.5652   # H122: 1999.5652 -> 1996.5652
$ws.Cells.Item(122, 9).Value = 2014   # I122: 2017.6316 -> 2014
$ws.Cells.Item(122, 11).Value = 6042   # K122: 6052.8948 -> 6042
$ws.Cells.Item(122, 13).Value = -3592   # M122: -3602.8948 -> -3592
$ws.Cells.Item(128, 8).Value = 80528.45   # H128: 80529.45 -> 80528.45
$ws.Cells.Item(128, 10).Value = 80528.45   # J128: 80529.45 -> 80528.45
$ws.Cells.Item(128, 12).Value = 80528.45   # L128: 80529.45 -> 80528.45
$ws.Cells.Item(128, 14).Value = -90488.45   # N128: -90489.45 -> -90488.45
$ws.Cells.Item(130, 8).Value = 122000   # H130: 62666.668 -> 122000
$ws.Cells.Item(130, 10).Value = 122000   # J130: 62666.668 -> 122000
$ws.Cells.Item(130, 12).Value = 122000   # L130: 62666.668 -> 122000
$ws.Cells.Item(130, 14).Value = -132040   # N130: -72706.66800000001 -> -132040
$ws.Cells.Item(132, 8).Value = 111116376   # H132: 142863200 -> 111116376
$ws.Cells.Item(132, 9).Value = 125005610   # I132: 142863200 -> 125005610
$ws.Cells.Item(132, 10).Value = 2470   # J132: 0 -> 2470
$ws.Cells.Item(132, 11).Value = 375016830   # K132: 428589600 -> 375016830
$ws.Cells.Item(132, 12).Value = 7410   # L132: 0 -> 7410
$ws.Cells.Item(132, 13).Value = -375014300   # M132: -428587070 -> -375014300
$ws.Cells.Item(132, 14).Value = -12470   # N132: None -> -12470

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 37464.3   # H7: 38685.344 -> 37464.3
$ws.Cells.Item(7, 9).Value = 3682.652   # I7: 3792.7144 -> 3682.652
$ws.Cells.Item(7, 10).Value = 148461.14   # J7: 130278.5 -> 148461.14
$ws.Cells.Item(7, 11).Value = 3682.652   # K7: 3792.7144 -> 3682.652
$ws.Cells.Item(7, 12).Value = 148461.14   # L7: 130278.5 -> 148461.14
$ws.Cells.Item(7, 13).Value = -3570.652   # M7: -3680.7144 -> -3570.652
$ws.Cells.Item(7, 14).Value = -148685.14   # N7: -130502.5 -> -148685.14
$ws.Cells.Item(40, 8).Value = 3280.25   # H40: 3176.1765 -> 3280.25
$ws.Cells.Item(40, 9).Value = 2652.6155   # I40: 2571.0715 -> 2652.6155
$ws.Cells.Item(40, 11).Value = 2652.6155   # K40: 2571.0715 -> 2652.6155
$ws.Cells.Item(40, 13).Value = -2516.6155   # M40: -2435.0715 -> -2516.6155
$ws.Cells.Item(122, 8).Value = 4868.9697   # H122: 4933.6875 -> 4868.9697
$ws.Cells.Item(122, 9).Value = 3988   # I122: 4033.7693 -> 3988
$ws.Cells.Item(122, 11).Value = 11964   # K122: 12101.3079 -> 11964
$ws.Cells.Item(122, 13).Value = -9514   # M122: -9651.3079 -> -9514
$ws.Cells.Item(126, 8).Value = 37464.3   # H126: 38685.344 -> 37464.3
$ws.Cells.Item(126, 9).Value = 3682.652   # I126: 3792.7144 -> 3682.652
$ws.Cells.Item(126, 10).Value = 148461.14   # J126: 130278.5 -> 148461.14
$ws.Cells.Item(126, 11).Value = 11047.956   # K126: 11378.1432 -> 11047.956
$ws.Cells.Item(126, 12).Value = 445383.42   # L126: 390835.5 -> 445383.42
$ws.Cells.Item(126, 13).Value = -8577.956   # M126: -8908.143199999999 -> -8577.956
$ws.Cells.Item(126, 14).Value = -450323.42   # N126: -395775.5 -> -450323.42
$ws.Cells.Item(132, 8).Value = 430742.28   # H132: 466483.4 -> 430742.28
$ws.Cells.Item(132, 9).Value = 528147   # I132: 557410.9 -> 528147
$ws.Cells.Item(132, 10).Value = 166358   # J132: 193701 -> 166358
$ws.Cells.Item(132, 11).Value = 1584441   # K132: 1672232.7 -> 1584441
$ws.Cells.Item(132, 12).Value = 499074   # L132: 581103 -> 499074
$ws.Cells.Item(132, 13).Value = -1581911   # M132: -1669702.7 -> -1581911
$ws.Cells.Item(132, 14).Value = -504134   # N132: -586163 -> -504134
$ws.Cells.Item(136, 8).Value = 60758.527   # H136: 55072 -> 60758.527
$ws.Cells.Item(136, 9).Value = 3039.2727   # I136: 2733.2307 -> 3039.2727
$ws.Cells.Item(136, 11).Value = 9117.8181   # K136: 8199.6921 -> 9117.8181
$ws.Cells.Item(136, 13).Value = -6567.8181   # M136: -5649.6921 -> -6567.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 65000   # H64: 0 -> 65000
$ws.Cells.Item(64, 10).Value = 65000   # J64: 0 -> 65000
$ws.Cells.Item(64, 12).Value = 65000   # L64: 0 -> 65000
$ws.Cells.Item(64, 14).Value = -65496   # N64: None -> -65496
$ws.Cells.Item(67, 8).Value = 65000   # H67: 0 -> 65000
$ws.Cells.Item(67, 10).Value = 65000   # J67: 0 -> 65000
$ws.Cells.Item(67, 12).Value = 65000   # L67: 0 -> 65000
$ws.Cells.Item(67, 14).Value = -66716   # N67: None -> -66716
$ws.Cells.Item(93, 8).Value = 105874   # H93: 100000 -> 105874
$ws.Cells.Item(93, 9).Value = 0   # I93: 90000 -> 0
$ws.Cells.Item(93, 10).Value = 105874   # J93: 110000 -> 105874
$ws.Cells.Item(93, 11).Value = 0   # K93: 90000 -> 0
$ws.Cells.Item(93, 12).Value = 105874   # L93: 110000 -> 105874
$ws.Cells.Item(93, 13).ClearContents()   # M93: remove (was -87504)
$ws.Cells.Item(93, 14).Value = -110866   # N93: -114992 -> -110866
$ws.Cells.Item(107, 8).Value = 23810980   # H107: 23810976 -> 23810980
$ws.Cells.Item(107, 9).Value = 35715870   # I107: 33334856 -> 35715870
$ws.Cells.Item(107, 10).Value = 1197   # J107: 1281.3334 -> 1197
$ws.Cells.Item(107, 11).Value = 107147610   # K107: 100004568 -> 107147610
$ws.Cells.Item(107, 12).Value = 3591   # L107: 3844.0002 -> 3591
$ws.Cells.Item(107, 13).Value = -107145690   # M107: -100002648 -> -107145690
$ws.Cells.Item(107, 14).Value = -7431   # N107: -7684.0002 -> -7431
$ws.Cells.Item(122, 8).Value = 5512.077   # H122: 5322.8887 -> 5512.077
$ws.Cells.Item(122, 9).Value = 2365.8667   # I122: 2243.25 -> 2365.8667
$ws.Cells.Item(122, 11).Value = 7097.6001   # K122: 6729.75 -> 7097.6001
$ws.Cells.Item(122, 13).Value = -4647.6001   # M122: -4279.75 -> -4647.6001
$ws.Cells.Item(124, 8).Value = 83954.5   # H124: 83979.75 -> 83954.5
$ws.Cells.Item(124, 10).Value = 83954.5   # J124: 83979.75 -> 83954.5
$ws.Cells.Item(124, 12).Value = 83954.5   # L124: 83979.75 -> 83954.5
$ws.Cells.Item(124, 14).Value = -93774.5   # N124: -93799.75 -> -93774.5
$ws.Cells.Item(136, 8).Value = 11548.042   # H136: 12021.869 -> 11548.042
$ws.Cells.Item(136, 9).Value = 1485.5625   # I136: 1541.2667 -> 1485.5625
$ws.Cells.Item(136, 11).Value = 4456.6875   # K136: 4623.800099999999 -> 4456.6875
$ws.Cells.Item(136, 13).Value = -1906.6875   # M136: -2073.800099999999 -> -1906.6875

Write-Host "Applied all corrections."